$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.837.63"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "2.906.97"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("D5").Value = "'590.14"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'144.73"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "2.905.34"
$ws.Range("E9").Value = "  -3.91%  "
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = "  -4.38%  "
$ws.Range("D12").Value = "'0.444"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").Value = "'33.42"
$ws.Range("E14").Value = "  -6.22%  "
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "3.390.24"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "60.794.68"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("E18").Value = "  -5.11%  "
$ws.Range("D19").Value = "2.906.52"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").Value = "'429.18"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "'7.08"
$ws.Range("E23").Value = "  -5.83%  "
$ws.Range("D24").Value = "'81.90"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'10.76"
$ws.Range("E25").Value = "  -6.26%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("D27").Value = "'11.99"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'2.27"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").Value = "'7.04"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").Value = "'26.49"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("D34").Value = "'0.108"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "0.0₃0852"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("D37").Value = "'5.60"
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("D39").Value = "'49.55"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("D41").Value = "'2.00"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").Value = "'8.62"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -4.76%  "
$ws.Range("D44").Value = "'40.06"
$ws.Range("E44").Value = "  -10.26%  "
$ws.Range("D45").Value = "'0.0349"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'372.78"
$ws.Range("E46").Value = "  -5.04%  "
$ws.Range("D47").Value = "2.698.87"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'131.33"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D50").Value = "'23.99"
$ws.Range("E50").Value = "  -9.22%  "
$ws.Range("D51").Value = "'0.107"
$ws.Range("E51").Value = "  -2.25%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
